$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.179.38"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "2.517.28"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  -2.01%  "

$ws.Range("D9").Value = "2.519.93"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("E11").Value = "  +1.66%  "

$ws.Range("E12").Value = "  -2.37%  "

$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("D14").Value = "2.962.90"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.81%  "

$ws.Range("D16").Value = "59.128.46"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").Value = "2.489.86"
$ws.Range("E18").Value = "  -0.69%  "

$ws.Range("E19").Value = "  -3.55%  "

$ws.Range("E20").Value = "  -1.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "319.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.417"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.57%  "

$ws.Range("E26").Value = "  +1.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.89%  "

$ws.Range("D30").Value = "0.0₃0763"
$ws.Range("E30").Value = "  -1.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.07%  "

$ws.Range("E33").Value = "  +0.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.75%  "

$ws.Range("E35").Value = "  -1.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.50%  "

$ws.Range("E38").Value = "  -3.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.25%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.800"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.596"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0925"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.53"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0506"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0222"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.52%  "
